# "Avance con el sistema Sobre Amortiguado"
#
# 1. Duplicate the "Implementación" sheet into a new "Hoja1" sheet that
#    explores the overdamped case (D1 = 2, pole product A20 = 5), keeping
#    only the left-hand (A:D) block.
# 2. On "Implementación": correct the C-column formula (drop the stray
#    factor of 2) and change the B9 capacitor value from 1000uF to 10uF.
# 3. Restore cursor/selection positions on every sheet, leaving
#    "Implementación" as the active tab.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Implementación")
$ws2 = $wb.Worksheets.Item("Polos Bremdow")

# --- 1. New "Hoja1" sheet: copy of "Implementación", overdamped variant ---
$ws1.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws3 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3.Name = "Hoja1"

# Only the A:D block is relevant on the new sheet; drop the copied
# "Polos Ruth" block that lived in F:I on the source sheet.
$ws3.Range("F1:I36").Clear() | Out-Null

# New R/C values that exercise the overdamped response.
$ws3.Range("D1").Value = 2
$ws3.Range("A20").Value = 5

# --- 2. Fix the C-column formula on "Implementación" (remove the /2) ---
$ws1.Range("C3:C4").Formula = "=3/(`$D`$1*B3)"
$ws1.Range("C5").Formula = "=3/(`$D`$1*B5)"
$ws1.Range("C6:C15").Formula = "=3/(`$D`$1*B6)"

# Capacitor in row 9 changed from 1000uF to 10uF.
$ws1.Range("B9").Formula = "=10*10^-6"

# --- 3. Selections, leaving "Implementación" as the active sheet/tab ---
$ws3.Range("F17").Select() | Out-Null
$ws2.Range("D30").Select() | Out-Null
$ws1.Activate()
$ws1.Range("E11").Select() | Out-Null
